$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 642.3125
$ws.Range("I15").Value = 642.3125
$ws.Range("K15").Value = 1926.9375
$ws.Range("M15").Value = -1757.9375
$ws.Range("H33").Value = 2646366.2
$ws.Range("I33").Value = 4116081
$ws.Range("K33").Value = 4116081
$ws.Range("M33").Value = -4115852
$ws.Range("H40").Value = 4488.778
$ws.Range("I40").Value = 3249
$ws.Range("J40").Value = 4587.96
$ws.Range("K40").Value = 3249
$ws.Range("L40").Value = 4587.96
$ws.Range("M40").Value = -3074
$ws.Range("N40").Value = -4937.96
$ws.Range("H43").Value = 71434584
$ws.Range("I43").Value = 250000660
$ws.Range("J43").Value = 8159.4
$ws.Range("K43").Value = 250000660
$ws.Range("L43").Value = 8159.4
$ws.Range("M43").Value = -250000591
$ws.Range("N43").Value = -8297.4
$ws.Range("H69").Value = 6509.846
$ws.Range("I69").Value = 4150
$ws.Range("J69").Value = 6817.6523
$ws.Range("K69").Value = 12450
$ws.Range("L69").Value = 20452.9569
$ws.Range("M69").Value = -11576
$ws.Range("N69").Value = -22200.9569
$ws.Range("H72").Value = 6509.846
$ws.Range("I72").Value = 4150
$ws.Range("J72").Value = 6817.6523
$ws.Range("K72").Value = 37350
$ws.Range("L72").Value = 61358.8707
$ws.Range("M72").Value = -32982
$ws.Range("N72").Value = -70094.8707
$ws.Range("H98").Value = 2132.147
$ws.Range("I98").Value = 2078
$ws.Range("J98").Value = 2384.8333
$ws.Range("K98").Value = 2078
$ws.Range("L98").Value = 2384.8333
$ws.Range("M98").Value = -580
$ws.Range("N98").Value = -5380.8333
$ws.Range("H118").Value = 83333790
$ws.Range("I118").Value = 111111464
$ws.Range("K118").Value = 333334392
$ws.Range("M118").Value = -333332735
$ws.Range("H122").Value = 2132.147
$ws.Range("I122").Value = 2078
$ws.Range("J122").Value = 2384.8333
$ws.Range("K122").Value = 6234
$ws.Range("L122").Value = 7154.499899999999
$ws.Range("M122").Value = -3784
$ws.Range("N122").Value = -12054.4999
$ws.Range("H132").Value = 31254270
$ws.Range("I132").Value = 38466100
$ws.Range("K132").Value = 115398300
$ws.Range("M132").Value = -115395770
$ws.Range("H137").Value = 51195.695
$ws.Range("I137").Value = 85941.05
$ws.Range("K137").Value = 257823.15
$ws.Range("M137").Value = -255273.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9878.091
$ws.Range("I32").Value = 6117.9785
$ws.Range("K32").Value = 6117.9785
$ws.Range("M32").Value = -5830.9785
$ws.Range("H45").Value = 5997800
$ws.Range("I45").Value = 10277274
$ws.Range("K45").Value = 10277274
$ws.Range("M45").Value = -10276897
$ws.Range("H61").Value = 6335.4585
$ws.Range("I61").Value = 6502.6816
$ws.Range("K61").Value = 6502.6816
$ws.Range("M61").Value = -6290.6816
$ws.Range("H63").Value = 5162.9414
$ws.Range("I63").Value = 1616.5
$ws.Range("K63").Value = 1616.5
$ws.Range("M63").Value = -930.5
$ws.Range("H66").Value = 5162.9414
$ws.Range("I66").Value = 1616.5
$ws.Range("K66").Value = 8082.5
$ws.Range("M66").Value = -4650.5
$ws.Range("H132").Value = 3010.6924
$ws.Range("J132").Value = 6065.6665
$ws.Range("L132").Value = 18196.9995
$ws.Range("N132").Value = -23256.9995
$ws.Range("H136").Value = 6335.4585
$ws.Range("I136").Value = 6502.6816
$ws.Range("K136").Value = 19508.0448
$ws.Range("M136").Value = -16958.0448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 596.1905
$ws.Range("I80").Value = 616.8570999999999
$ws.Range("J80").Value = 585.8570999999999
$ws.Range("K80").Value = 616.8570999999999
$ws.Range("L80").Value = 585.8570999999999
$ws.Range("M80").Value = 381.1429000000001
$ws.Range("N80").Value = -2581.8571
$ws.Range("H83").Value = 596.1905
$ws.Range("I83").Value = 616.8570999999999
$ws.Range("J83").Value = 585.8570999999999
$ws.Range("K83").Value = 3084.2855
$ws.Range("L83").Value = 2929.2855
$ws.Range("M83").Value = 1907.7145
$ws.Range("N83").Value = -12913.2855
$ws.Range("H107").Value = 6501367.5
$ws.Range("I107").Value = 7944560.5
$ws.Range("K107").Value = 7944560.5
$ws.Range("M107").Value = -7942640.5
$ws.Range("H134").Value = 3609.111
$ws.Range("J134").Value = 10644.857
$ws.Range("L134").Value = 31934.571
$ws.Range("N134").Value = -37004.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3810.9473
$ws.Range("I99").Value = 3552.8572
$ws.Range("K99").Value = 3552.8572
$ws.Range("M99").Value = -2054.8572
$ws.Range("H105").Value = 1133.9375
$ws.Range("I105").Value = 762
$ws.Range("K105").Value = 762
$ws.Range("M105").Value = 985
$ws.Range("H126").Value = 3810.9473
$ws.Range("I126").Value = 3552.8572
$ws.Range("K126").Value = 10658.5716
$ws.Range("M126").Value = -8188.571599999999
$ws.Range("H132").Value = 97145.7
$ws.Range("I132").Value = 85089.414
$ws.Range("K132").Value = 255268.242
$ws.Range("M132").Value = -252738.242
$ws.Range("H138").Value = 89999
$ws.Range("J138").Value = 89999
$ws.Range("L138").Value = 89999
$ws.Range("N138").Value = -100279
$ws.Range("H141").Value = 348884
$ws.Range("J141").Value = 348884
$ws.Range("L141").Value = 348884
$ws.Range("N141").Value = -359244

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 42
$ws.Range("I6").Value = 42
$ws.Range("K6").Value = 126
$ws.Range("M6").Value = -13
$ws.Range("H45").Value = 98.888885
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 98.888885
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 296.666655
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -1360.666655
$ws.Range("H59").Value = 998.5
$ws.Range("I59").Value = 998.5
$ws.Range("K59").Value = 2995.5
$ws.Range("M59").Value = -2455.5
$ws.Range("H64").Value = 733.3333
$ws.Range("J64").Value = 900
$ws.Range("L64").Value = 2700
$ws.Range("N64").Value = -3240
$ws.Range("H67").Value = 733.3333
$ws.Range("J67").Value = 900
$ws.Range("L67").Value = 2700
$ws.Range("N67").Value = -4572
$ws.Range("H128").Value = 199990.25
$ws.Range("I128").Value = 199990.25
$ws.Range("K128").Value = 599970.75
$ws.Range("M128").Value = -594990.75
$ws.Range("H129").Value = 2223240.2
$ws.Range("I129").Value = 2857818.8
$ws.Range("J129").Value = 2216
$ws.Range("K129").Value = 8573456.399999999
$ws.Range("L129").Value = 6648
$ws.Range("M129").Value = -8568456.399999999
$ws.Range("N129").Value = -16648
$ws.Range("H137").Value = 2027.125
$ws.Range("I137").Value = 1495.3846
$ws.Range("J137").Value = 4331.3335
$ws.Range("K137").Value = 4486.1538
$ws.Range("L137").Value = 12994.0005
$ws.Range("M137").Value = 613.8462
$ws.Range("N137").Value = -23194.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3480.348
$ws.Range("I132").Value = 3411.5
$ws.Range("K132").Value = 10234.5
$ws.Range("M132").Value = -7704.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 37455.04
$ws.Range("I22").Value = 60165.066
$ws.Range("J22").Value = 3390
$ws.Range("K22").Value = 60165.066
$ws.Range("L22").Value = 3390
$ws.Range("M22").Value = -59870.066
$ws.Range("N22").Value = -3980
$ws.Range("H27").Value = 37455.04
$ws.Range("I27").Value = 60165.066
$ws.Range("J27").Value = 3390
$ws.Range("K27").Value = 60165.066
$ws.Range("L27").Value = 3390
$ws.Range("M27").Value = -60058.066
$ws.Range("N27").Value = -3604
$ws.Range("H40").Value = 3199.739
$ws.Range("I40").Value = 2258.8235
$ws.Range("J40").Value = 5865.6665
$ws.Range("K40").Value = 2258.8235
$ws.Range("L40").Value = 5865.6665
$ws.Range("M40").Value = -2122.8235
$ws.Range("N40").Value = -6137.6665
$ws.Range("H61").Value = 13889799
$ws.Range("I61").Value = 15873842
$ws.Range("K61").Value = 15873842
$ws.Range("M61").Value = -15873640
$ws.Range("H100").Value = 2327.111
$ws.Range("I100").Value = 1932.4706
$ws.Range("J100").Value = 2998
$ws.Range("K100").Value = 1932.4706
$ws.Range("L100").Value = 2998
$ws.Range("M100").Value = -1391.4706
$ws.Range("N100").Value = -4080
$ws.Range("H113").Value = 13889799
$ws.Range("I113").Value = 15873842
$ws.Range("K113").Value = 15873842
$ws.Range("M113").Value = -15871672

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12825027
$ws.Range("I81").Value = 23811008
$ws.Range("J81").Value = 8050
$ws.Range("K81").Value = 47622016
$ws.Range("L81").Value = 16100
$ws.Range("M81").Value = -47620955
$ws.Range("N81").Value = -18222
$ws.Range("H84").Value = 12825027
$ws.Range("I84").Value = 23811008
$ws.Range("J84").Value = 8050
$ws.Range("K84").Value = 238110080
$ws.Range("L84").Value = 80500
$ws.Range("M84").Value = -238104776
$ws.Range("N84").Value = -91108
$ws.Range("H126").Value = 2346.4546
$ws.Range("I126").Value = 2346.4546
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7039.3638
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4569.3638
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 3205.125
$ws.Range("I136").Value = 3150.4827
$ws.Range("J136").Value = 3733.3333
$ws.Range("K136").Value = 9451.4481
$ws.Range("L136").Value = 11199.9999
$ws.Range("M136").Value = -6901.4481
$ws.Range("N136").Value = -16299.9999
